# Butterfly matrix.xlsx edit script
# Commit: "rename generations cols in Butterfly matrix, add Poaceae column to
#          forage_plant_locations spreadsheet"
#
# For this workbook, the visible edits are:
#   1) Rename several header labels on "Full Species Set" (shorten/standardise
#      the "Adult/Larval reliance" + "Adult/Larval generations" + the
#      "*_damp plant" headers to the Ad_/L_ naming convention used elsewhere).
#   2) Two data corrections: Brimstone's adult-generations count, and Marsh
#      fritillary's larval-generations count (both become 2).
#   3) A brand new "larvae forage plants" sheet added after "Key", seeded with
#      the Species column (header + the 24 species names) ready to be filled
#      in - companion to the "forage_plant_locations" work mentioned in the
#      commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Full Species Set")

# --- 1) Header label renames (row 4) -------------------------------------
$ws.Range("B4").Value = "Ad_reliance"
$ws.Range("C4").Value = "L_reliance"
$ws.Range("D4").Value = "Ad_generations"
$ws.Range("E4").Value = "L_generations"

$ws.Range("Z4").Value  = "Ad_H_damp"
$ws.Range("AA4").Value = "Ad_M_damp"
$ws.Range("AB4").Value = "Ad_CA_damp"
$ws.Range("AC4").Value = "Ad_CG_damp"
$ws.Range("AD4").Value = "Ad_R_damp"

$ws.Range("AY4").Value = "L_H_damp"
$ws.Range("AZ4").Value = "L_M_damp"
$ws.Range("BA4").Value = "L_CA_damp"
$ws.Range("BB4").Value = "L_CG_damp"
$ws.Range("BC4").Value = "L_R_damp"

# --- 2) Data corrections ---------------------------------------------------
# Brimstone (row 5): Ad_generations 1 -> 2
$ws.Range("D5").Value = 2
# Marsh fritillary (row 16): L_generations "1 (2 according to supp?)" -> 2
$ws.Range("E16").Value = 2

# --- 3) Add the new "larvae forage plants" sheet --------------------------
$keySheet = $wb.Worksheets.Item("Key")
$newSheet = $wb.Worksheets.Add($null, $keySheet)
$newSheet.Name = "larvae forage plants"

# Column A: Species header + the same 24 species as "Full Species Set"
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newSheet.Range("A1").Value = $ws.Range("A4").Value2

$destRow = 1
for ($r = 5; $r -le $lastRow; $r++) {
    $destRow = $destRow + 1
    $newSheet.Cells.Item($destRow, 1).Value = $ws.Cells.Item($r, 1).Value2
}

$newSheet.Columns.Item(1).ColumnWidth = 18.33203125
$newSheet.Range("A1:A" + $destRow).HorizontalAlignment = -4108
$newSheet.Range("A1").Select()

# --- Restore the original sheet/selection focus ----------------------------
$ws.Activate()
$ws.Range("AG19").Select()
